$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (after "calendar_year") shifting everything else down
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "total_smolt"
$ws.Range("B3").Value = "The total number of smolt counted at the juvenile weir at Auke Creek"

$ws.Range("B3").Select()
